$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.420.21"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.572.90"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "291.27"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.3771"
$ws.Range("E7").Value = "  +2.96%  "
$ws.Range("D8").Value = "49.81"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "0.07663"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "6.003"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "6.929"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "1.569.77"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "90.26"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "0.06763"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "16.80"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").Value = "6.222"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "2.428"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("D25").Value = "22.429.75"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "2.735"
$ws.Range("E26").Value = "  -8.45%  "
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").Value = "146.05"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "5.037"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").Value = "126.38"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").Value = "1.745.43"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "6.207"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").Value = "1.001"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").Value = "10.04"
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").Value = "0.08577"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("D37").Value = "0.02546"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "0.2320"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.345"
$ws.Range("E39").Value = "  +7.52%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06578"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("D41").Value = "5.463"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").Value = "0.6451"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "14.13"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "0.6021"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").Value = "1.306"
$ws.Range("E48").Value = "  +8.66%  "
$ws.Range("D49").Value = "2.086"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").Value = "125.27"
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("D51").Value = "0.07331"
$ws.Range("E51").Value = "  +0.81%  "
